$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: locate the paragraph containing a (unique) piece of text and
# return its full Range, so we are robust against paragraph-index shifts
# caused by earlier edits in this same script.
# ---------------------------------------------------------------------------
function Get-ParaRangeByText($text) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) { throw "Text not found: $text" }
    return $rng.Paragraphs(1).Range
}

$pkgHeader = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-ParaXml($range, $bodyInnerXml) {
    $xml = $pkgHeader + '<w:body>' + $bodyInnerXml + '</w:body>' + $pkgFooter
    $range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 1) "By: Itamar" -> "By: " + proofErr-wrapped "Itamar"
# ---------------------------------------------------------------------------
$r = Get-ParaRangeByText("By: Itamar")
$body = '<w:p><w:r><w:t xml:space="preserve">By: </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Itamar</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/></w:p>'
Set-ParaXml $r $body

# ---------------------------------------------------------------------------
# 2) "...(Brazil is huge and the weather probably differs a lot)." paragraph
#    -> wrap "huge" in gramStart/gramEnd proofErr markers.
# ---------------------------------------------------------------------------
$r = Get-ParaRangeByText("In order for the weather to be relevant")
$body = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">In order for the weather to be relevant we need to pull the weather separately for each row because they are in different city' + [char]0x2019 + 's (Brazil is </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>huge</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> and the weather probably differs a lot).</w:t></w:r></w:p>'
Set-ParaXml $r $body

# ---------------------------------------------------------------------------
# 3) "...will be a waist of memory..." paragraph
#    -> wrap "waist" in spellStart/spellEnd proofErr markers.
# ---------------------------------------------------------------------------
$r = Get-ParaRangeByText("Create new columns in our own data set")
$body = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Create new columns in our own data set and insert the weather values in. However, this will be a </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>waist</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> of memory and will cost us in the time or each RUN</w:t></w:r></w:p>'
Set-ParaXml $r $body

# ---------------------------------------------------------------------------
# 4) "...like the guy from github did)." paragraph
#    -> wrap "github" in spellStart/spellEnd proofErr markers.
# ---------------------------------------------------------------------------
$r = Get-ParaRangeByText("The better thing to do is to create the columns")
$body = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">The better thing to do is to create the columns in our dataset but keep them empty and to pull/point towards the data in the weather data. This way we will have only 60 rows of weather data instead of having 110000 rows of duplicated weather data if we used option 1 (like the guy from </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>github</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> did).</w:t></w:r></w:p>'
Set-ParaXml $r $body

# ---------------------------------------------------------------------------
# 5) Structural change around the "Re-write the coodbook" item:
#    - A new (empty) green-highlighted paragraph holding the "_GoBack"
#      bookmark is inserted right after "I don't know how to do this...".
#    - The "Re-write the coodbook (...)" paragraph follows it, now with a
#      plain "ind left=360" pPr (no longer a numbered ListParagraph item)
#      and "coodbook" wrapped in spellStart/spellEnd proofErr markers.
#    - The paragraph's old slot (right before "Add columns by date...") is
#      removed (we fold it into the two new paragraphs above).
# ---------------------------------------------------------------------------
$r = Get-ParaRangeByText("Re-write the coodbook")
$body = '<w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' + `
    '<w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">Re-write the </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>coodbook</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:hint="cs"/><w:highlight w:val="green"/><w:rtl/><w:lang w:bidi="he-IL"/></w:rPr><w:t>בוצע בירוק</w:t></w:r>' + `
    '<w:r><w:rPr><w:highlight w:val="green"/><w:lang w:bidi="he-IL"/></w:rPr><w:t>)</w:t></w:r></w:p>'
Set-ParaXml $r $body

# ---------------------------------------------------------------------------
# 6) Last paragraph ("Ckv ckv" + the old "_GoBack" bookmark location) is
#    emptied out -- the text and the bookmark both moved away (the
#    bookmark now lives on the new paragraph inserted in step 5).
#    NOTE: this is the very last paragraph in the document body, and
#    InsertXML-ing its whole Range has proven to leave a stray extra
#    empty paragraph behind (an engine quirk around the final paragraph
#    mark), so it is edited with plain text Find/Replace + the Bookmarks
#    collection instead of XML surgery.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Ckv ckv", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

Write-Output "All steps done"
